$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("D1").Value = "watched"
$ws.Range("E1").Value = "rate"
$ws.Range("F1").Value = "review"

# Copy the header style (bold, border, alignment) from C1 to the new header cells
$ws.Range("C1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update row 2 (Interstellar entry)
$ws.Range("B2").Value = "Interstellar"
$ws.Range("C2").Value = "20xx"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "★★★★★★★★★★ (10.0)"
$ws.Range("F2").Value = "-"

# Remove row 3 (Oppenheimer) entirely, shifting rows up
$ws.Rows("3:3").Delete()
